$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "refri"
$ws.Range("A3").Value = "doritos"
